$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1641.1045
$ws.Range("J17").Value = 1726.6809
$ws.Range("L17").Value = 5180.0427
$ws.Range("N17").Value = -5516.0427
$ws.Range("H19").Value = 889007.6
$ws.Range("I19").Value = 1212197.8
$ws.Range("J19").Value = 234.75
$ws.Range("K19").Value = 1212197.8
$ws.Range("L19").Value = 234.75
$ws.Range("M19").Value = -1212022.8
$ws.Range("N19").Value = -584.75
$ws.Range("H40").Value = 1465.8889
$ws.Range("I40").Value = 2001
$ws.Range("J40").Value = 1399
$ws.Range("K40").Value = 2001
$ws.Range("L40").Value = 1399
$ws.Range("M40").Value = -1826
$ws.Range("N40").Value = -1749
$ws.Range("H64").Value = 2830.6924
$ws.Range("I64").Value = 2819.9
$ws.Range("J64").Value = 2866.6667
$ws.Range("K64").Value = 2819.9
$ws.Range("L64").Value = 2866.6667
$ws.Range("M64").Value = -2571.9
$ws.Range("N64").Value = -3362.6667
$ws.Range("H67").Value = 2830.6924
$ws.Range("I67").Value = 2819.9
$ws.Range("J67").Value = 2866.6667
$ws.Range("K67").Value = 2819.9
$ws.Range("L67").Value = 2866.6667
$ws.Range("M67").Value = -1961.9
$ws.Range("N67").Value = -4582.6667
$ws.Range("H76").Value = 3600
$ws.Range("I76").Value = 3600
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 3600
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = -3285
$ws.Range("N76").ClearContents()
$ws.Range("H79").Value = 3600
$ws.Range("I79").Value = 3600
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 3600
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = -2508
$ws.Range("N79").ClearContents()
$ws.Range("H107").Value = 949.5625
$ws.Range("I107").Value = 945.8
$ws.Range("K107").Value = 945.8
$ws.Range("M107").Value = 974.2
$ws.Range("H135").Value = 1276
$ws.Range("I135").Value = 1195.6666
$ws.Range("J135").Value = 1396.5
$ws.Range("K135").Value = 10760.9994
$ws.Range("L135").Value = 12568.5
$ws.Range("M135").Value = -8225.999400000001
$ws.Range("N135").Value = -17638.5
$ws.Range("H137").Value = 4225.8975
$ws.Range("I137").Value = 4066.56
$ws.Range("J137").Value = 4510.4287
$ws.Range("K137").Value = 12199.68
$ws.Range("L137").Value = 13531.2861
$ws.Range("M137").Value = -9649.68
$ws.Range("N137").Value = -18631.2861
$ws.Range("H138").Value = 2374.7058
$ws.Range("I138").Value = 1255.0883
$ws.Range("J138").Value = 3494.3235
$ws.Range("K138").Value = 3765.2649
$ws.Range("L138").Value = 10482.9705
$ws.Range("M138").Value = 1374.7351
$ws.Range("N138").Value = -20762.9705

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 410.76923
$ws.Range("I2").Value = 384.87097
$ws.Range("J2").Value = 511.125
$ws.Range("K2").Value = 384.87097
$ws.Range("L2").Value = 511.125
$ws.Range("M2").Value = -271.87097
$ws.Range("N2").Value = -737.125
$ws.Range("H32").Value = 4674.62
$ws.Range("I32").Value = 2508.7568
$ws.Range("J32").Value = 10839
$ws.Range("K32").Value = 2508.7568
$ws.Range("L32").Value = 10839
$ws.Range("M32").Value = -2221.7568
$ws.Range("N32").Value = -11413
$ws.Range("H116").Value = 410.76923
$ws.Range("I116").Value = 384.87097
$ws.Range("J116").Value = 511.125
$ws.Range("K116").Value = 384.87097
$ws.Range("L116").Value = 511.125
$ws.Range("M116").Value = 1909.12903
$ws.Range("N116").Value = -5099.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 410.76923
$ws.Range("I3").Value = 384.87097
$ws.Range("J3").Value = 511.125
$ws.Range("K3").Value = 384.87097
$ws.Range("L3").Value = 511.125
$ws.Range("M3").Value = -270.87097
$ws.Range("N3").Value = -739.125
$ws.Range("H105").Value = 2607.7
$ws.Range("I105").Value = 2592.5789
$ws.Range("J105").Value = 2895
$ws.Range("K105").Value = 2592.5789
$ws.Range("L105").Value = 2895
$ws.Range("M105").Value = -845.5789
$ws.Range("N105").Value = -6389

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2747.1177
$ws.Range("I31").Value = 1143.6897
$ws.Range("J31").Value = 4860.727
$ws.Range("K31").Value = 1143.6897
$ws.Range("L31").Value = 4860.727
$ws.Range("M31").Value = -848.6896999999999
$ws.Range("N31").Value = -5450.727
$ws.Range("H34").Value = 2747.1177
$ws.Range("I34").Value = 1143.6897
$ws.Range("J34").Value = 4860.727
$ws.Range("K34").Value = 1143.6897
$ws.Range("L34").Value = 4860.727
$ws.Range("M34").Value = -941.6896999999999
$ws.Range("N34").Value = -5264.727
$ws.Range("H62").Value = 62503250
$ws.Range("I62").Value = 62503250
$ws.Range("K62").Value = 62503250
$ws.Range("M62").Value = -62502626
$ws.Range("H65").Value = 62503250
$ws.Range("I65").Value = 62503250
$ws.Range("K65").Value = 312516250
$ws.Range("M65").Value = -312513130
$ws.Range("H86").Value = 1988.6842
$ws.Range("J86").Value = 2291.3635
$ws.Range("L86").Value = 2291.3635
$ws.Range("N86").Value = -4537.363499999999
$ws.Range("H89").Value = 1988.6842
$ws.Range("J89").Value = 2291.3635
$ws.Range("L89").Value = 11456.8175
$ws.Range("N89").Value = -22688.8175
$ws.Range("H99").Value = 4350418.5
$ws.Range("I99").Value = 7694148.5
$ws.Range("J99").Value = 3569.3
$ws.Range("K99").Value = 7694148.5
$ws.Range("L99").Value = 3569.3
$ws.Range("M99").Value = -7692650.5
$ws.Range("N99").Value = -6565.3
$ws.Range("H126").Value = 4350418.5
$ws.Range("I126").Value = 7694148.5
$ws.Range("J126").Value = 3569.3
$ws.Range("K126").Value = 23082445.5
$ws.Range("L126").Value = 10707.9
$ws.Range("M126").Value = -23079975.5
$ws.Range("N126").Value = -15647.9
$ws.Range("H132").Value = 2540.5122
$ws.Range("I132").Value = 1274.3182
$ws.Range("J132").Value = 4006.6316
$ws.Range("K132").Value = 3822.9546
$ws.Range("L132").Value = 12019.8948
$ws.Range("M132").Value = -1292.9546
$ws.Range("N132").Value = -17079.8948

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 746.1
$ws.Range("I113").Value = 655.8889
$ws.Range("K113").Value = 1967.6667
$ws.Range("M113").Value = 202.3332999999998
$ws.Range("H123").Value = 4524.75
$ws.Range("I123").Value = 4366.6665
$ws.Range("J123").Value = 4999
$ws.Range("K123").Value = 13099.9995
$ws.Range("L123").Value = 14997
$ws.Range("M123").Value = -10649.9995
$ws.Range("N123").Value = -19897
$ws.Range("H134").Value = 4261.96
$ws.Range("I134").Value = 3767.7856
$ws.Range("J134").Value = 4890.909
$ws.Range("K134").Value = 11303.3568
$ws.Range("L134").Value = 14672.727
$ws.Range("M134").Value = -6233.356800000001
$ws.Range("N134").Value = -24812.727

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 3988
$ws.Range("I70").Value = 3760
$ws.Range("K70").Value = 3760
$ws.Range("M70").Value = -3490
$ws.Range("H73").Value = 3988
$ws.Range("I73").Value = 3760
$ws.Range("K73").Value = 3760
$ws.Range("M73").Value = -2824
$ws.Range("H80").Value = 20835992
$ws.Range("I80").Value = 31252438
$ws.Range("K80").Value = 31252438
$ws.Range("M80").Value = -31251440
$ws.Range("H83").Value = 20835992
$ws.Range("I83").Value = 31252438
$ws.Range("K83").Value = 156262190
$ws.Range("M83").Value = -156257198
$ws.Range("H102").Value = 3065.1482
$ws.Range("I102").Value = 2366.3125
$ws.Range("K102").Value = 2366.3125
$ws.Range("M102").Value = -744.3125

